{"js": "// Append a new paragraph (item \"6.\") after the last paragraph of the\n// document body (the \"csv file\" paragraph), per the commit\n// \"update in reading csv and input\".\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst newText =\n  \"6.in validation methods I used generic functions that can be used widely in the system. For example in order to validate priority enum I created a versatile function that takes in a raw value represented as a string and compares it against possible values of an enum class. This is better than to create a function that only checks if a value is within the range of the priority enum because this will allow us to reuse this function to match against other enum classes and will help us achieving a good SE practice of not repeating our selves.\";\n\nlastParagraph.insertParagraph(newText, \"After\");\n\nawait context.sync();\n", "ps1": "# Append a new paragraph (item \"6.\") after the last paragraph of the\n# document body (the \"csv file\" paragraph), per the commit\n# \"update in reading csv and input\".\n\n$d = $word.ActiveDocument\n\n$newText = \"6.in validation methods I used generic functions that can be used widely in the system. For example in order to validate priority enum I created a versatile function that takes in a raw value represented as a string and compares it against possible values of an enum class. This is better than to create a function that only checks if a value is within the range of the priority enum because this will allow us to reuse this function to match against other enum classes and will help us achieving a good SE practice of not repeating our selves.\"\n\n$rng = $d.Content\n$rng.Collapse(0)  # wdCollapseEnd\n$rng.InsertParagraphAfter()\n$rng.Collapse(0)\n\n$newPara = $d.Paragraphs.Last\n$newPara.Range.Text = $newText\n"}
